# "I will go to home from university"
# B10 (quantity for ЕКГ-12 / row 28) drops from 8 to 5 and row 22's quantity
# (which used to be hard-typed as 8) now follows it via a formula, so every
# downstream total / triangle-load computation recalculates automatically.
# A brand-new row 24 is also added that folds the auxiliary-load total (D23)
# into the overall apparent load (H15) and derives a transformer-sizing
# figure from it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- core edit -----------------------------------------------------------
# Quantity of ЕКГ-12 excavators: 8 -> 5. All shared formulas in D10, H10,
# I10, D14, H14, I14, H15, K15, L15 recompute automatically.
$ws.Range("B10").Value = 5

# Row 22 ("Трансформатор ... ЕКГ-12.5") quantity now tracks B10 instead of
# being a fixed literal.
$ws.Range("B22").Formula = "=B10"

# --- new row 24: fold the auxiliary loads (D23) into the apparent load
# (H15) and size the transformer from the combined figure.
$ws.Range("G24").Formula = "=D23+H15"
$ws.Range("H24").Formula = "=G24/(2*0.8)"

# --- view state: scroll down and move the active selection ---------------
$ws.Range("B11").Select()
$excel.ActiveWindow.ScrollRow = 20
$excel.ActiveWindow.ScrollColumn = 1
